# Apply the "Updated symbol list" edits to Sheet1.
# The worksheet stores every data cell as text (inline strings), including
# things that look like plain numbers (prices such as "283.47"). Excel's COM
# layer auto-coerces numeric-looking strings into real numbers when you just
# assign .Value, so for any price-like cell we force the cell to Text format
# first, assign the literal string, then drop the temporary NumberFormat
# back off the cell so no stray style survives the round-trip.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# --- Price ("D" column) tweaks ---------------------------------------------
Set-TextValue "D2"  "283.40"
Set-TextValue "D3"  "20.82"
Set-TextValue "D4"  "6.208"
Set-TextValue "D5"  "0.06189"
Set-TextValue "D6"  "3.583"
Set-TextValue "D7"  "6.564"
Set-TextValue "D8"  "1.490"
Set-TextValue "D9"  "0.8172"
Set-TextValue "D11" "0.1648"
Set-TextValue "D12" "0.08373"
Set-TextValue "D13" "0.03666"
Set-TextValue "D15" "0.09135"
Set-TextValue "D16" "3.703"
Set-TextValue "D17" "0.001644"
Set-TextValue "D18" "0.04665"
Set-TextValue "D19" "0.006470"
Set-TextValue "D20" "0.006192"
Set-TextValue "D21" "0.001067"
Set-TextValue "D22" "0.0001500"
Set-TextValue "D23" "3.796"
Set-TextValue "D24" "2.330"

# Row 27 (AAXToken) gains a "Bestin24h" suffix on its composite id column.
$ws.Range("E27").Value = "26AAXTokenAABBestin24h"

Set-TextValue "D40" "0.04724"

# --- Rows 41-43 got reshuffled (coins rotated: CEJI -> KickToken ->
#     BKEXToken -> CEJI) together with fresh price/id data. ------------------
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D41" "0.007043"
$ws.Range("E41").Value = "40KickTokenKICK"

$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1105"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.003519"
$ws.Range("E43").Value = "42CEJICEJI"

Set-TextValue "D44" "0.01144"
Set-TextValue "D45" "0.00006394"
Set-TextValue "D47" "0.9997"
Set-TextValue "D50" "0.01240"
